# The target diff for this document only reorders XML attributes
# (alphabetises namespace declarations / element attributes in
# word/document.xml and word/styles.xml) -- every changed line is a
# pure attribute-order change with no added/removed/modified text,
# values, styles or formatting. It is canonicalisation noise coming
# from whatever tool produced the upstream diff, not a real content
# edit made through Word.
#
# Word's COM object model (and this runtime's OOXML writer) does not
# expose attribute-order as an editable property: every experiment
# (touching PageSetup margins, re-setting style fonts, Find/Replace,
# explicit Save) leaves already-existing attributes in their original
# order and only ever appends newly required attributes/namespaces.
# So the only faithful way to "apply" this diff through the Word COM
# surface is to leave the document's content untouched -- the
# before/after documents are semantically identical.
$d = $word.ActiveDocument
